# Daily attendance processing - 2026-01-05 22:34:51
# Normalize the "Recorded By" (column G) lists so that entries are ordered
# by a fixed priority (System accounts first, then the reviewer/backup
# accounts, then individual user emails last), while keeping the relative
# order stable for any ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priority = @{
    "System" = 0;
    "system" = 1;
    "admin@admin.com" = 2;
    "backup@backdoor.com" = 3;
    "dnasr281@gmail.com" = 4;
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text
    if ([string]::IsNullOrEmpty($text)) { continue }

    $parts = $text -split ", "
    if ($parts.Count -le 1) { continue }

    $indexed = @()
    for ($i = 0; $i -lt $parts.Count; $i++) {
        $part = $parts[$i]
        $pri = $priority[$part]
        if ($null -eq $pri) { $pri = 99 }
        $indexed += [PSCustomObject]@{ Part = $part; Pri = $pri }
    }

    $sorted = $indexed | Sort-Object -Property @("Pri")
    $newParts = @()
    foreach ($s in $sorted) {
        $newParts += $s.Part
    }
    $newText = $newParts -join ", "

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
